# Insert one new row of data just above the current row 669
# (2026/12/29 ...) so it becomes row 670, etc. The new row holds a
# 2026/01/19 / 月 / 1 / 201 entry, matching the existing 2026/01/19
# entries already present in rows 666-668.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 669 (and everything below it) down by one row.
$ws.Rows("669").Insert()

# Column A holds a date formatted as plain text (e.g. "2026/12/29"),
# not a real date serial. Force the cell to Text format before writing
# so Excel doesn't auto-convert the literal string into a date value,
# then drop back to the default "Normal" style so no stray formatting
# is left behind on the cell.
$ws.Range("A669").NumberFormat = "@"
$ws.Range("A669").Value = "2026/01/19"
$ws.Range("A669").Style = "Normal"

$ws.Range("B669").Value = "月"
$ws.Range("C669").Value = 1
$ws.Range("D669").Value = 201
